$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add note to existing row 21 (D21)
$ws.Range("D21").Value = "MF working, similar mse, need to format new dataframe to train and test with"

# Add new row 22
$ws.Range("A21").Copy()
$ws.Range("A22").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A22").Value = 45436
$ws.Range("B22").Value = 5
$ws.Range("D22").Value = "Finished mf using recosystem, not very good.. Looking to optimize"

$ws.Range("D22").Select()

# Row height updates (text-wrap reflow side effect of the edit)
$ws.Rows.Item(9).RowHeight = 30
$ws.Rows.Item(11).RowHeight = 30
$ws.Rows.Item(12).RowHeight = 45
$ws.Rows.Item(13).RowHeight = 30
$ws.Rows.Item(14).RowHeight = 45
$ws.Rows.Item(17).RowHeight = 30
$ws.Rows.Item(18).RowHeight = 45
$ws.Rows.Item(19).RowHeight = 30
$ws.Rows.Item(20).RowHeight = 45
$ws.Rows.Item(21).RowHeight = 30
$ws.Rows.Item(22).RowHeight = 30
